# Retrospective.pptx - "Added my parts to Readme and PPT"
#
# Fills in Joel's two placeholder bullet points ("Joel-") on the
# "Individual Portions" slides with his actual retrospective text.
#
# Slide 2 ("What is impeding us?") gets his "impediment" writeup and
# slide 4 ("What can we do to improve?") gets his "improvement" writeup.
# Both paragraphs are located by matching their existing "Joel-" text so
# the script is resilient to minor shape/paragraph reshuffling.

$p = $ppt.ActivePresentation

function Set-JoelParagraph {
    param(
        [int]$slideIndex,
        [string]$newText
    )

    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item("Content Placeholder 2")
    $tr = $shape.TextFrame.TextRange

    $paraCount = $tr.Paragraphs().Count
    for ($i = 1; $i -le $paraCount; $i++) {
        $para = $tr.Paragraphs($i, 1)
        $existing = $para.Text.TrimEnd("`r")
        if ($existing -eq "Joel-") {
            $run = $para.Runs(1, 1)
            $run.Text = $newText
            return
        }
    }
}

# Slide 2 ("Individual Portions" - What is impeding us?)
Set-JoelParagraph 2 "Joel- Research what resources are available online to help my coding"

# Slide 4 ("Individual Portions" - What can we do to improve?)
Set-JoelParagraph 4 "Joel- I need to start looking at other Twitch bot API's so that I can see what goals I can achieve in regards to the functionality of the Twitch bot code."
